# Append a new "8^9" equation row (row 6) below the existing data (rows 1-5),
# matching the layout/text-typing of the prior rows exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 6

# Leading apostrophe forces each value to be stored as literal text (not a
# number), matching the existing rows where every cell is a plain string
# ("8^9", "1.342177e+08", "1648538283957" all stay textual rather than
# getting parsed as a formula / float / integer).
$ws.Cells.Item($newRow, 1).Value = "'8^9"
$ws.Cells.Item($newRow, 2).Value = "'1.342177e+08"
$ws.Cells.Item($newRow, 3).Value = "'1648538283957"

# The apostrophe entry leaves a "quote prefix" style on the cells; copy the
# plain (unstyled) look from the row above so the new row matches the
# formatting of the other data rows.
$ws.Cells.Item($newRow, 1).Style = $ws.Cells.Item($newRow - 1, 1).Style
$ws.Cells.Item($newRow, 2).Style = $ws.Cells.Item($newRow - 1, 2).Style
$ws.Cells.Item($newRow, 3).Style = $ws.Cells.Item($newRow - 1, 3).Style
